$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to remain text even when the new value looks numeric
    # (e.g. "122.80"), matching the workbook's inlineStr/text convention,
    # then restore the original ("Normal") cell style so no formatting
    # side effects leak into the saved file.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 (Bitcoin)
Set-TextValue "D2" "61.861.89"
$ws.Range("E2").Value = "  +8.42%  "

# Row 3 (Ethereum)
Set-TextValue "D3" "3.436.74"
$ws.Range("E3").Value = "  +5.39%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 (BNB)
Set-TextValue "D5" "413.25"
$ws.Range("E5").Value = "  +4.04%  "

# Row 6 (Solana)
Set-TextValue "D6" "122.80"
$ws.Range("E6").Value = "  +13.21%  "

# Row 7 (LidoStakedEther)
Set-TextValue "D7" "3.433.72"
$ws.Range("E7").Value = "  +5.45%  "

# Row 8 (XRP)
Set-TextValue "D8" "0.589"
$ws.Range("E8").Value = "  +1.41%  "

# Row 9 (USDC)
Set-TextValue "D9" "0.999"
$ws.Range("E9").Value = "  -0.04%  "

# Row 10 (Cardano)
Set-TextValue "D10" "0.649"
$ws.Range("E10").Value = "  +4.47%  "

# Row 11 (Dogecoin)
Set-TextValue "D11" "0.125"
$ws.Range("E11").Value = "  +31.02%  "

# Row 12 (Avalanche)
Set-TextValue "D12" "41.35"
$ws.Range("E12").Value = "  +4.72%  "

# Row 13 (TRON)
$ws.Range("E13").Value = "  -0.29%  "

# Row 14 (WrappedliquidstakedEther2.0)
Set-TextValue "D14" "3.962.83"
$ws.Range("E14").Value = "  +4.88%  "

# Row 15 (Polkadot)
Set-TextValue "D15" "8.48"
$ws.Range("E15").Value = "  +2.29%  "

# Row 16 (Chainlink)
Set-TextValue "D16" "19.60"
$ws.Range("E16").Value = "  +3.15%  "

# Row 17 (WrappedEther)
Set-TextValue "D17" "3.412.60"
$ws.Range("E17").Value = "  +4.64%  "

# Row 18 (WrappedBTC)
Set-TextValue "D18" "61.713.85"
$ws.Range("E18").Value = "  +8.40%  "

# Row 19 (Polygon)
$ws.Range("E19").Value = "  -0.45%  "

# Row 20 (Uniswap)
Set-TextValue "D20" "10.76"
$ws.Range("E20").Value = "  -2.49%  "

# Row 21 (ShibaInu)
Set-TextValue "D21" "0.0000133"
$ws.Range("E21").Value = "  +22.85%  "

# Row 22 (ImmutableX)
$ws.Range("E22").Value = "  -0.95%  "

# Row 23 (BitcoinCash)
Set-TextValue "D23" "314.67"
$ws.Range("E23").Value = "  +7.39%  "

# Row 24 - swapped with row 25: now Litecoin
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "81.35"
$ws.Range("E24").Value = "  +9.40%  "

# Row 25 - swapped with row 24: now InternetComputer(DFINITY)
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D25" "12.95"
$ws.Range("E25").Value = "  +0.13%  "

# Row 26 (PancakeSwap)
$ws.Range("E26").Value = "  -0.68%  "

# Row 27 (EthereumClassic)
Set-TextValue "D27" "30.99"
$ws.Range("E27").Value = "  +10.03%  "

# Row 28 (RenderToken)
Set-TextValue "D28" "7.93"
$ws.Range("E28").Value = "  +6.96%  "

# Row 29 (LEO)
Set-TextValue "D29" "4.29"
$ws.Range("E29").Value = "  -1.96%  "

# Row 30 (Kaspa)
$ws.Range("E30").Value = "  +2.68%  "

# Row 31 (Filecoin)
Set-TextValue "D31" "7.69"
$ws.Range("E31").Value = "  -3.30%  "

# Row 32 (Hedera)
$ws.Range("E32").Value = "  +4.48%  "

# Row 33 (Toncoin)
Set-TextValue "D33" "2.56"
$ws.Range("E33").Value = "  +20.15%  "

# Row 34 (InjectiveProtocol)
Set-TextValue "D34" "42.02"
$ws.Range("E34").Value = "  +5.31%  "

# Row 35 (Cosmos)
Set-TextValue "D35" "11.37"
$ws.Range("E35").Value = "  +1.68%  "

# Row 36 (Dai)
$ws.Range("E36").Value = "  +0.11%  "

# Row 37 (VeChain)
Set-TextValue "D37" "0.0479"
$ws.Range("E37").Value = "  -1.35%  "

# Row 38 (OKB)
Set-TextValue "D38" "52.44"
$ws.Range("E38").Value = "  +2.17%  "

# Row 39 (LidoDAOToken)
Set-TextValue "D39" "3.51"
$ws.Range("E39").Value = "  +1.47%  "

# Row 40 (FirstDigitalUSD)
$ws.Range("E40").Value = "  -0.24%  "

# Row 41 (Stacks)
Set-TextValue "D41" "3.01"
$ws.Range("E41").Value = "  -0.77%  "

# Row 42 (ARBITRUM)
Set-TextValue "D42" "1.98"
$ws.Range("E42").Value = "  +5.77%  "

# Row 43 (Stellar)
$ws.Range("E43").Value = "  +2.61%  "

# Row 44 (Monero)
Set-TextValue "D44" "134.30"
$ws.Range("E44").Value = "  -2.17%  "

# Row 45 (Celestia)
Set-TextValue "D45" "17.17"
$ws.Range("E45").Value = "  +2.56%  "

# Row 46 (TheGraph)
Set-TextValue "D46" "0.283"
$ws.Range("E46").Value = "  -0.68%  "

# Row 47 (NEARProtocol)
Set-TextValue "D47" "3.87"
$ws.Range("E47").Value = "  -1.19%  "

# Row 48 - swapped with row 49: now WEMIXToken
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D48" "2.20"
$ws.Range("E48").Value = "  -0.65%  "

# Row 49 - swapped with row 48: now EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "21.81"
$ws.Range("E49").Value = "  -2.26%  "

# Row 50 (Maker)
Set-TextValue "D50" "2.206.46"
$ws.Range("E50").Value = "  +2.49%  "

# Row 51 (ApeXProtocol)
$ws.Range("E51").Value = "  +3.48%  "
